$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global Parameters")

# Update the Power Scaling Factor (row 12) and Cost Scaling Factor (row 15)
# values to 1, reflecting the new "investment case study" scaling.
$ws.Range("C12").Value = 1
$ws.Range("C15").Value = 1

# Update the active selection to C16, as left by the author after editing.
$ws.Range("C16").Select()
